$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")

# Clear out the existing header columns D:G (vital_status, id, name, description)
$ws.Range("D1:G1").ClearContents()

# Remove the data validation on column D
$ws.Range("D2:D1048576").Validation.Delete()

# Set the new header row values
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "description"
